$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 15000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 15000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 15000
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -15460
$ws.Range("H33").Value = 138.21428
$ws.Range("I33").Value = 149
$ws.Range("J33").Value = 111.25
$ws.Range("K33").Value = 149
$ws.Range("L33").Value = 111.25
$ws.Range("M33").Value = 80
$ws.Range("N33").Value = -569.25
$ws.Range("H70").Value = 3148.7778
$ws.Range("J70").Value = 3262.7144
$ws.Range("L70").Value = 9788.143199999999
$ws.Range("N70").Value = -10328.1432
$ws.Range("H73").Value = 3148.7778
$ws.Range("J73").Value = 3262.7144
$ws.Range("L73").Value = 9788.143199999999
$ws.Range("N73").Value = -11660.1432
$ws.Range("H87").Value = 116499
$ws.Range("J87").Value = 158998
$ws.Range("L87").Value = 158998
$ws.Range("N87").Value = -161494
$ws.Range("H90").Value = 116499
$ws.Range("J90").Value = 158998
$ws.Range("L90").Value = 476994
$ws.Range("N90").Value = -489474
$ws.Range("H111").Value = 11021.5625
$ws.Range("I111").Value = 13738.3
$ws.Range("K111").Value = 41214.89999999999
$ws.Range("M111").Value = -38147.89999999999

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 280
$ws.Range("I4").Value = 280
$ws.Range("K4").Value = 280
$ws.Range("M4").Value = -164
$ws.Range("H32").Value = 4056220.2
$ws.Range("I32").Value = 3892677
$ws.Range("K32").Value = 3892677
$ws.Range("M32").Value = -3892390
$ws.Range("H35").Value = 39768.5
$ws.Range("I35").Value = 19691.334
$ws.Range("K35").Value = 19691.334
$ws.Range("M35").Value = -19285.334
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H122").Value = 3307.4443
$ws.Range("J122").Value = 3523.75
$ws.Range("L122").Value = 10571.25
$ws.Range("N122").Value = -15471.25

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 892
$ws.Range("I20").Value = 1019.25
$ws.Range("J20").Value = 637.5
$ws.Range("K20").Value = 1019.25
$ws.Range("L20").Value = 637.5
$ws.Range("M20").Value = -772.25
$ws.Range("N20").Value = -1131.5
$ws.Range("H22").Value = 383.66666
$ws.Range("I22").Value = 475.5
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 475.5
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -302.5
$ws.Range("N22").Value = -546
$ws.Range("H25").Value = 400
$ws.Range("I25").Value = 400
$ws.Range("K25").Value = 400
$ws.Range("M25").Value = -165
$ws.Range("H75").Value = 67107
$ws.Range("I75").Value = 67107
$ws.Range("K75").Value = 67107
$ws.Range("M75").Value = -66171
$ws.Range("H78").Value = 67107
$ws.Range("I78").Value = 67107
$ws.Range("K78").Value = 201321
$ws.Range("M78").Value = -196641
$ws.Range("H81").Value = 41955.2
$ws.Range("J81").Value = 41955.2
$ws.Range("L81").Value = 41955.2
$ws.Range("N81").Value = -44077.2
$ws.Range("H84").Value = 41955.2
$ws.Range("J84").Value = 41955.2
$ws.Range("L84").Value = 125865.6
$ws.Range("N84").Value = -136473.6
$ws.Range("H99").Value = 2224.5557
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 3673.6667
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 3673.6667
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -6669.6667
$ws.Range("H102").Value = 64500
$ws.Range("I102").Value = 64500
$ws.Range("K102").Value = 64500
$ws.Range("M102").Value = -61255

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2598.2856
$ws.Range("I2").Value = 899.6667
$ws.Range("J2").Value = 3872.25
$ws.Range("K2").Value = 899.6667
$ws.Range("L2").Value = 3872.25
$ws.Range("M2").Value = -786.6667
$ws.Range("N2").Value = -4098.25
$ws.Range("H31").Value = 2440.5557
$ws.Range("I31").Value = 1995
$ws.Range("K31").Value = 1995
$ws.Range("M31").Value = -1700
$ws.Range("H34").Value = 2440.5557
$ws.Range("I34").Value = 1995
$ws.Range("K34").Value = 1995
$ws.Range("M34").Value = -1793
$ws.Range("H58").Value = 2841
$ws.Range("I58").Value = 1991
$ws.Range("J58").Value = 3266
$ws.Range("K58").Value = 1991
$ws.Range("L58").Value = 3266
$ws.Range("M58").Value = -1788
$ws.Range("N58").Value = -3672
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H74").Value = 44750
$ws.Range("J74").Value = 44750
$ws.Range("L74").Value = 44750
$ws.Range("N74").Value = -46498
$ws.Range("H77").Value = 44750
$ws.Range("J77").Value = 44750
$ws.Range("L77").Value = 134250
$ws.Range("N77").Value = -142986
$ws.Range("H136").Value = 2841
$ws.Range("I136").Value = 1991
$ws.Range("J136").Value = 3266
$ws.Range("K136").Value = 5973
$ws.Range("L136").Value = 9798
$ws.Range("M136").Value = -3423
$ws.Range("N136").Value = -14898

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 144.25
$ws.Range("J2").Value = 163.42857
$ws.Range("L2").Value = 980.57142
$ws.Range("N2").Value = -1206.57142
$ws.Range("H16").Value = 250
$ws.Range("I16").Value = 250
$ws.Range("K16").Value = 750
$ws.Range("M16").Value = -577
$ws.Range("H120").Value = 5555
$ws.Range("I120").Value = 5555
$ws.Range("K120").Value = 16665
$ws.Range("M120").Value = -11827
$ws.Range("H140").Value = 1178
$ws.Range("I140").Value = 1178
$ws.Range("K140").Value = 3534
$ws.Range("M140").Value = 1646

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 53467.285
$ws.Range("J15").Value = 53467.285
$ws.Range("L15").Value = 53467.285
$ws.Range("N15").Value = -54043.285
$ws.Range("H70").Value = 4999.5
$ws.Range("I70").Value = 4999.5
$ws.Range("K70").Value = 4999.5
$ws.Range("M70").Value = -4729.5
$ws.Range("H73").Value = 4999.5
$ws.Range("I73").Value = 4999.5
$ws.Range("K73").Value = 4999.5
$ws.Range("M73").Value = -4063.5
$ws.Range("H81").Value = 53467.285
$ws.Range("J81").Value = 53467.285
$ws.Range("L81").Value = 53467.285
$ws.Range("N81").Value = -55463.285
$ws.Range("H84").Value = 53467.285
$ws.Range("J84").Value = 53467.285
$ws.Range("L84").Value = 160401.855
$ws.Range("N84").Value = -170385.855
$ws.Range("H107").Value = 1649.8182
$ws.Range("I107").Value = 414.14285
$ws.Range("J107").Value = 3812.25
$ws.Range("K107").Value = 414.14285
$ws.Range("L107").Value = 3812.25
$ws.Range("M107").Value = 1505.85715
$ws.Range("N107").Value = -7652.25
$ws.Range("H122").Value = 3945.875
$ws.Range("I122").Value = 4802.8335
$ws.Range("K122").Value = 14408.5005
$ws.Range("M122").Value = -11958.5005

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 865.125
$ws.Range("I61").Value = 1150.75
$ws.Range("J61").Value = 579.5
$ws.Range("K61").Value = 1150.75
$ws.Range("L61").Value = 579.5
$ws.Range("M61").Value = -948.75
$ws.Range("N61").Value = -983.5
$ws.Range("H93").Value = 3699.5
$ws.Range("I93").Value = 3500
$ws.Range("J93").Value = 3899
$ws.Range("K93").Value = 3500
$ws.Range("L93").Value = 3899
$ws.Range("M93").Value = -2252
$ws.Range("N93").Value = -6395
$ws.Range("H94").Value = 20000
$ws.Range("J94").Value = 20000
$ws.Range("L94").Value = 20000
$ws.Range("N94").Value = -21352
$ws.Range("H100").Value = 1992.125
$ws.Range("I100").Value = 1862.5714
$ws.Range("J100").Value = 2899
$ws.Range("K100").Value = 1862.5714
$ws.Range("L100").Value = 2899
$ws.Range("M100").Value = -1321.5714
$ws.Range("N100").Value = -3981
$ws.Range("H113").Value = 865.125
$ws.Range("I113").Value = 1150.75
$ws.Range("J113").Value = 579.5
$ws.Range("K113").Value = 1150.75
$ws.Range("L113").Value = 579.5
$ws.Range("M113").Value = 1019.25
$ws.Range("N113").Value = -4919.5

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H74").Value = 21693.143
$ws.Range("I74").Value = 20227.5
$ws.Range("K74").Value = 20227.5
$ws.Range("M74").Value = -19291.5
$ws.Range("H77").Value = 21693.143
$ws.Range("I77").Value = 20227.5
$ws.Range("K77").Value = 60682.5
$ws.Range("M77").Value = -56002.5
$ws.Range("H107").Value = 1849.4375
$ws.Range("I107").Value = 1727.9286
$ws.Range("K107").Value = 5183.7858
$ws.Range("M107").Value = -3263.7858
$ws.Range("H126").Value = 3299.5
$ws.Range("I126").Value = 2699.25
$ws.Range("K126").Value = 8097.75
$ws.Range("M126").Value = -5627.75
